$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns I1 ("I0") and J1 ("IF"), matching the formatting
# (bold, centered, bordered) of the existing header cells by copying the
# style from the last existing header cell (H1).
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Corresponding data-row values.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
